$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 73, shifting existing rows 73:99 down to 74:100
$ws.Rows.Item(73).Insert()

# Fill in the new row 73 with the inserted data
$ws.Cells.Item(73, 1).Value = 2
$ws.Cells.Item(73, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(73, 3).Value = "Coquimbo"
$ws.Cells.Item(73, 4).Value = 44559
$ws.Cells.Item(73, 5).Value = 4
$ws.Cells.Item(73, 6).Value = "Fruta"
$ws.Cells.Item(73, 7).Value = 100109
$ws.Cells.Item(73, 8).Value = "Uva"
$ws.Cells.Item(73, 9).Value = 100109001
$ws.Cells.Item(73, 10).Value = "Uva"
$ws.Cells.Item(73, 11).Value = "Flame Seedless"
$ws.Cells.Item(73, 12).Value = "Primera"
$ws.Cells.Item(73, 13).Value = 900
$ws.Cells.Item(73, 14).Value = 6000
$ws.Cells.Item(73, 15).Value = 7000
$ws.Cells.Item(73, 16).Value = 6500
$ws.Cells.Item(73, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(73, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(73, 19).Value = 650
$ws.Cells.Item(73, 20).Value = 10
